# Apply the change recorded in the commit "Add files via upload":
# a new value "A15289126" (the student ID) is written into cell C4 of
# Sheet1, which becomes a new shared-string entry (index 21).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C4").Value = "A15289126"

# Keep the view/selection in sync with the saved file (the diff shows the
# active cell moving from B6 to C4, and the top-left cell moving to A2).
$excel.Goto($ws.Range("A2"), $false)
$ws.Range("C4").Select()
